# Adds two new columns, I (header "I0") and J (header "IF"), to Sheet1.
# Header cells (row 1) get the same bold/centered/bordered style as the
# existing header cells (B1:H1); data cells (rows 2:66) are plain numbers,
# matching the formatting of the existing data columns (e.g. column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, centered, bordered -> style index 1)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data values (rows 2-66) -------------------------------------------
$I0 = @(12,5,6,7,7,6,8,7,9,5,8,6,5,8,7,8,6,8,7,8,8,5,8,7,8,7,7,6,7,9,8,9,8,6,6,9,9,5,7,7,5,4,9,8,7,7,9,5,7,8,10,7,9,9,9,8,8,7,7,4,9,5,4,6,6)
$IF = @(13,5,6,7,7,6,8,7,9,5,8,7,5,8,7,8,7,8,7,8,8,6,8,7,8,7,7,6,7,9,8,9,8,6,7,9,9,5,8,7,6,4,9,8,7,8,9,6,7,8,10,7,9,9,9,8,8,7,7,5,9,5,5,6,6)

for ($i = 0; $i -lt $I0.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I0[$i]
    $ws.Cells.Item($row, 10).Value = $IF[$i]
}
